$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03835166666666667
$ws.Range("H2").Value = 0.115055
$ws.Range("I2").Value = 0.0181239951898282
$ws.Range("J2").Value = 0.0181239951898282
$ws.Range("M2").Value = 24.91851366666667
$ws.Range("N2").Value = 74.75554099999999
$ws.Range("O2").Value = 0.2924799159147552
$ws.Range("P2").Value = 0.2924799159147553
$ws.Range("Q2").Value = 0.9556665299727777
$ws.Range("R2").Value = 8.600998769755
$ws.Range("S2").Value = 0.005300904589160378
$ws.Range("T2").Value = 0.005300904589160379
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03835166666666667
$ws.Range("H3").Value = 0.115055
$ws.Range("I3").Value = 0.0181239951898282
$ws.Range("J3").Value = 0.0181239951898282
$ws.Range("O3").Value = 0.4753125595076708
$ws.Range("P3").Value = 0.4753125595076708
$ws.Range("Q3").Value = 1.553064944567222
$ws.Range("R3").Value = 13.977584501105
$ws.Range("S3").Value = 0.008614562542181953
$ws.Range("T3").Value = 0.008614562542181953
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03835166666666667
$ws.Range("H4").Value = 0.115055
$ws.Range("I4").Value = 0.0181239951898282
$ws.Range("J4").Value = 0.0181239951898282
$ws.Range("M4").Value = 19.78346566666667
$ws.Range("N4").Value = 59.350397
$ws.Range("O4").Value = 0.232207524577574
$ws.Range("P4").Value = 0.232207524577574
$ws.Range("Q4").Value = 0.7587288807594444
$ws.Range("R4").Value = 6.828559926835
$ws.Range("S4").Value = 0.004208528058485864
$ws.Range("T4").Value = 0.004208528058485864
$ws.Range("I5").Value = 0.3727881574250648
$ws.Range("J5").Value = 0.3727881574250648
$ws.Range("M5").Value = 24.91851366666667
$ws.Range("N5").Value = 74.75554099999999
$ws.Range("O5").Value = 0.2924799159147552
$ws.Range("P5").Value = 0.2924799159147553
$ws.Range("Q5").Value = 19.65687813806655
$ws.Range("R5").Value = 176.911903242599
$ws.Range("S5").Value = 0.1090330489376995
$ws.Range("T5").Value = 0.1090330489376995
$ws.Range("I6").Value = 0.3727881574250648
$ws.Range("J6").Value = 0.3727881574250648
$ws.Range("O6").Value = 0.4753125595076708
$ws.Range("P6").Value = 0.4753125595076708
$ws.Range("S6").Value = 0.177190893259856
$ws.Range("T6").Value = 0.177190893259856
$ws.Range("I7").Value = 0.3727881574250648
$ws.Range("J7").Value = 0.3727881574250648
$ws.Range("M7").Value = 19.78346566666667
$ws.Range("N7").Value = 59.350397
$ws.Range("O7").Value = 0.232207524577574
$ws.Range("P7").Value = 0.232207524577574
$ws.Range("Q7").Value = 15.60611435177589
$ws.Range("R7").Value = 140.455029165983
$ws.Range("S7").Value = 0.08656421522750926
$ws.Range("T7").Value = 0.08656421522750928
$ws.Range("G8").Value = 1.288873333333333
$ws.Range("H8").Value = 3.86662
$ws.Range("I8").Value = 0.6090878473851071
$ws.Range("J8").Value = 0.609087847385107
$ws.Range("M8").Value = 24.91851366666667
$ws.Range("N8").Value = 74.75554099999999
$ws.Range("O8").Value = 0.2924799159147552
$ws.Range("P8").Value = 0.2924799159147553
$ws.Range("Q8").Value = 32.11680777126889
$ws.Range("R8").Value = 289.0512699414199
$ws.Range("S8").Value = 0.1781459623878954
$ws.Range("T8").Value = 0.1781459623878954
$ws.Range("G9").Value = 1.288873333333333
$ws.Range("H9").Value = 3.86662
$ws.Range("I9").Value = 0.6090878473851071
$ws.Range("J9").Value = 0.609087847385107
$ws.Range("O9").Value = 0.4753125595076708
$ws.Range("P9").Value = 0.4753125595076708
$ws.Range("Q9").Value = 52.19340294609111
$ws.Range("R9").Value = 469.74062651482
$ws.Range("S9").Value = 0.2895071037056328
$ws.Range("T9").Value = 0.2895071037056328
$ws.Range("G10").Value = 1.288873333333333
$ws.Range("H10").Value = 3.86662
$ws.Range("I10").Value = 0.6090878473851071
$ws.Range("J10").Value = 0.609087847385107
$ws.Range("M10").Value = 19.78346566666667
$ws.Range("N10").Value = 59.350397
$ws.Range("O10").Value = 0.232207524577574
$ws.Range("P10").Value = 0.232207524577574
$ws.Range("Q10").Value = 25.49838133868223
$ws.Range("R10").Value = 229.48543204814
$ws.Range("S10").Value = 0.1414347812915789
$ws.Range("T10").Value = 0.1414347812915789

Write-Host "Applied new TPM-based values to rows 2-10"
